# Add new "online resource" attribute rows to the reader configuration /
# default-value lookup sheets, and rename the postal-code keys to match
# the updated reader schema.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Config_AB / Config_BC: rename contact postal code attribute keys
# ---------------------------------------------------------------------
$cfgAB = $wb.Worksheets.Item("Config_AB")
$cfgAB.Range("A86").Value2 = "contact_postal_code"
$cfgAB.Range("B86").Value2 = "contacts{}.postal_code"

$cfgBC = $wb.Worksheets.Item("Config_BC")
$cfgBC.Range("A79").Value2 = "contact_postal_code"
$cfgBC.Range("B79").Value2 = "contacts{}.postal_code"

# ---------------------------------------------------------------------
# Default_AB: add the GeoDiscover Alberta online-resource attribute block
# ---------------------------------------------------------------------
$defAB = $wb.Worksheets.Item("Default_AB")
$defAB.Range("A40").Value2 = "online_resource_link"
$defAB.Range("B40").Value2 = "https://geodiscover.alberta.ca/geoportal/"
$defAB.Range("A41").Value2 = "online_resource_protocol"
$defAB.Range("B41").Value2 = "HTTPS"
$defAB.Range("A42").Value2 = "online_resource_description"
$defAB.Range("B42").Value2 = "GeoDiscover Alberta provides enhanced details regarding Alberta's geospatial data."
$defAB.Range("A43").Value2 = "online_resource_description_other_lang_locale"
$defAB.Range("B43").Value2 = "#fra"
$defAB.Range("A44").Value2 = "online_resource_description_other_lang"
$defAB.Range("B44").Value2 = "GéoDécouvrez l'Alberta fournit des détails améliorés sur les données géospatiales de l'Alberta."

$geoLink = $defAB.Hyperlinks.Add($defAB.Range("B40"), "https://geodiscover.alberta.ca/geoportal/")
$defAB.Range("B40").Style = "Hyperlink"

# ---------------------------------------------------------------------
# Default_BC: add the DataBC online-resource attribute block
# ---------------------------------------------------------------------
$defBC = $wb.Worksheets.Item("Default_BC")
$defBC.Range("A24").Value2 = "online_resource_link"
$defBC.Range("B24").Value2 = "https://data.gov.bc.ca/"
$defBC.Range("A25").Value2 = "online_resource_protocol"
$defBC.Range("B25").Value2 = "HTTPS"
$defBC.Range("A26").Value2 = "online_resource_description"
$defBC.Range("B26").Value2 = "DataBC encourages and enables the strategic management and sharing of data across the government enterprise and with the public. "
$defBC.Range("A27").Value2 = "online_resource_description_other_lang_locale"
$defBC.Range("B27").Value2 = "#fra"
$defBC.Range("A28").Value2 = "online_resource_description_other_lang"
$defBC.Range("B28").Value2 = "DataBC encourage et permet la gestion et le partage stratégiques des données dans l'ensemble de l'entreprise gouvernementale et avec le public. "

# ---------------------------------------------------------------------
# Restore view/selection state to match the saved workbook
# ---------------------------------------------------------------------
$defBC.Activate()
$defBC.Range("C27").Select()

$defAB.Activate()
$defAB.Range("A40:A44").Select()

$cfgAB.Activate()
$cfgAB.Range("A86").Select()

# Config_BC is the sheet that was active/selected when the workbook was
# last saved, so leave it activated last.
$cfgBC.Activate()
$cfgBC.Range("B79").Select()
